# New weekly price record inserted for "Ají" at Feria Lagunitas de Puerto Montt.
# This pushes the existing data rows (320-423) down by one (to 321-424) and
# populates the newly opened row 320 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 320, shifting rows 320:423 down to 321:424.
$ws.Rows("320:320").Insert()

# Populate the new row 320 with the latest weekly record.
$ws.Cells.Item(320, 1).Value  = 4
$ws.Cells.Item(320, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(320, 3).Value  = "Los Lagos"
$ws.Cells.Item(320, 4).Value  = 45093
$ws.Cells.Item(320, 5).Value  = 10
$ws.Cells.Item(320, 6).Value  = 100112021
$ws.Cells.Item(320, 7).Value  = "Ají"
$ws.Cells.Item(320, 8).Value  = "Inferno"
$ws.Cells.Item(320, 9).Value  = "Primera"
$ws.Cells.Item(320, 10).Value = 160
$ws.Cells.Item(320, 11).Value = 19000
$ws.Cells.Item(320, 12).Value = 20000
$ws.Cells.Item(320, 13).Value = 19500
$ws.Cells.Item(320, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(320, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(320, 16).Value = 1950
$ws.Cells.Item(320, 17).Value = 10
$ws.Cells.Item(320, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same numeric date style used by the rest of
# column D (style index referenced by D321, the row just below).
$ws.Cells.Item(320, 4).NumberFormat = $ws.Cells.Item(321, 4).NumberFormat
